$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bus voltage magnitude results for the 380 kV case (vm_pu.xlsx)
# Columns order matches B,C,D,E,F,I,J,K,L,M,N (column G stays 1, H stays blank)
$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

$data = @{
    2 = @(1.02, 1.019615537826321, 1.030244161831632, 1.02073141874446, 1.031869944366426, 1.032330490091489, 1.024817299578962, 1.033055735294763, 1.023570837009088, 1.034676814831682, 1.012352941027119)
    3 = @(1.02, 1.020418741547745, 1.030694025018843, 1.021408986067245, 1.032926554786676, 1.03246672452827, 1.025258099855024, 1.033314876543066, 1.024055059310531, 1.035541415050812, 1.012500233965003)
    4 = @(1.02, 1.020939073971036, 1.030985466045683, 1.021848334882879, 1.033611166214083, 1.032553921726561, 1.025543283017283, 1.033482178482455, 1.024368622493415, 1.036101187076476, 1.012595486470674)
    5 = @(1.02, 1.021157965242964, 1.031108069715108, 1.022033255049133, 1.033899194559104, 1.03259035000628, 1.025663162364387, 1.033552420173679, 1.024500500423476, 1.036336590285594, 1.012635516834974)
    6 = @(1.02, 1.021194726422295, 1.031128660158618, 1.022064316687998, 1.033947568520677, 1.032596453002365, 1.025683289918909, 1.033564208632449, 1.024522646554056, 1.036376119890786, 1.01264223728868)
    7 = @(1.02, 1.020941998244901, 1.030987103962531, 1.021850804941056, 1.033615014008843, 1.032554409385304, 1.025544884898156, 1.033483117418458, 1.024370384434417, 1.036104332252091, 1.01259602141306)
    8 = @(1.02, 1.019886857904411, 1.030396121762669, 1.020960214798792, 1.032226840943467, 1.032376728626752, 1.024966278565028, 1.033143391252603, 1.023734431974528, 1.034968944074454, 1.012402730686257)
    9 = @(1.02, 1.018032274231474, 1.029357491519772, 1.019397983538111, 1.029787748963793, 1.032056342581662, 1.023946411150804, 1.032541891369779, 1.022615687678885, 1.032970725972796, 1.01206171802566)
    10 = @(1.02, 1.016799142070423, 1.028667031705742, 1.018361373525942, 1.028166488002801, 1.031837889509978, 1.0232663705159, 1.032139036343794, 1.021871199677545, 1.03164031026892, 1.011834125335449)
    11 = @(1.02, 1.016265973901493, 1.028368542739083, 1.017913688436188, 1.027465615542395, 1.031742151822812, 1.022971888560029, 1.031964169866265, 1.021549161954295, 1.03106464852868, 1.011735520840771)
    12 = @(1.02, 1.016068051154956, 1.028257745324816, 1.017747576379306, 1.027205453200642, 1.031706419046319, 1.022862502964752, 1.031899153576529, 1.021429593716022, 1.030850885761259, 1.011698886838305)
    13 = @(1.02, 1.016110500816851, 1.028281508337171, 1.017783199890988, 1.027261251066228, 1.031714091600184, 1.022885966613888, 1.031913102628136, 1.021455239185008, 1.030896735705145, 1.011706745307164)
    14 = @(1.02, 1.016249611095803, 1.028359382656121, 1.01789995391804, 1.027444106909294, 1.031739201635915, 1.02296284674469, 1.031958796879952, 1.021539277360584, 1.03104697752784, 1.01173249282152)
    15 = @(1.02, 1.016335337430696, 1.028407373545258, 1.017971913582381, 1.027556793369339, 1.031754650046736, 1.023010214913775, 1.031986942290273, 1.021591062837954, 1.031139554934943, 1.011748355680075)
    16 = @(1.02, 1.016834543380928, 1.02868685180612, 1.01839110979328, 1.028213026791618, 1.031844219213863, 1.023285913988705, 1.032150632745082, 1.021892579331161, 1.031678523912676, 1.011840668260047)
    17 = @(1.02, 1.017147893653193, 1.028862292215279, 1.018654376024186, 1.028624971719452, 1.03190009720793, 1.023458848190914, 1.032253197774623, 1.022081802032406, 1.032016717255478, 1.011898558974198)
    18 = @(1.02, 1.017330741318874, 1.028964670278932, 1.018808047907976, 1.028865362561679, 1.031932579274369, 1.023559715793945, 1.032312980887833, 1.022192204236256, 1.032214019882796, 1.011932320279345)
    19 = @(1.02, 1.017393100450115, 1.028999586407438, 1.018860465187086, 1.02894734831146, 1.031943636036521, 1.023594108668678, 1.032333358338727, 1.022229853890878, 1.032281301764173, 1.011943831088651)
    20 = @(1.02, 1.0171142662879, 1.028843464281993, 1.018626118333307, 1.028580762507624, 1.031894113465491, 1.023440294185672, 1.032242197784636, 1.022061496928672, 1.031980428148123, 1.011892348400791)
    21 = @(1.02, 1.016208643270606, 1.028336448533781, 1.01786556783606, 1.027390255640316, 1.031731812086299, 1.022940207505553, 1.0319453427942, 1.021514528796988, 1.031002733291228, 1.011724911035114)
    22 = @(1.02, 1.015639935159, 1.028018100589021, 1.01738841095861, 1.026642737541165, 1.031628774629408, 1.022625772799821, 1.0317583338023, 1.021170923399884, 1.030388386387004, 1.011619590999372)
    23 = @(1.02, 1.015941351885744, 1.028186821197661, 1.017641262496624, 1.027038915831171, 1.031683490561918, 1.022792461201269, 1.031857504973311, 1.021353046723103, 1.03071402803199, 1.011675427319494)
    24 = @(1.02, 1.017129460812919, 1.028851971670743, 1.018638886418103, 1.028600738406621, 1.031896817603015, 1.023448677948277, 1.032247168334486, 1.022070671832914, 1.031996825509626, 1.011895154708223)
    25 = @(1.02, 1.018511161197496, 1.029625664648847, 1.019801005552053, 1.030417471231893, 1.032140029932876, 1.024210098851081, 1.032697725385257, 1.022904678858383, 1.033487012912515, 1.012149923993779)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $values[$i]
    }
}
